# Update the "想去人数" (want-to-go count) figures in column F for both the
# "展览" sheet and the "全部类型" sheet, reflecting the latest scrape values.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 531
$wsExhibit.Range("F5").Value = 253
$wsExhibit.Range("F7").Value = 238
$wsExhibit.Range("F8").Value = 2279
$wsExhibit.Range("F9").Value = 384
$wsExhibit.Range("F10").Value = 5667

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 531
$wsAll.Range("F6").Value = 253
$wsAll.Range("F8").Value = 238
$wsAll.Range("F11").Value = 2279
$wsAll.Range("F12").Value = 384
$wsAll.Range("F13").Value = 5667
